$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.818.04'
$ws.Range("E2").Value = '  +1.06%  '

# Row 3
$ws.Range("D3").Value = '3.199.75'
$ws.Range("E3").Value = '  +0.62%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = '599.25'
$ws.Range("E5").Value = '  +4.30%  '

# Row 6
$ws.Range("D6").Value = '151.05'
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").Value = '3.197.62'
$ws.Range("E8").Value = '  +0.59%  '

# Row 9
$ws.Range("E9").Value = '  +1.67%  '

# Row 10
$ws.Range("E10").Value = '  -1.34%  '

# Row 11
$ws.Range("D11").Value = '6.15'
$ws.Range("E11").Value = '  -1.14%  '

# Row 12
$ws.Range("D12").Value = '0.510'
$ws.Range("E12").Value = '  +0.57%  '

# Row 13
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  -0.81%  '

# Row 14
$ws.Range("D14").Value = '38.38'
$ws.Range("E14").Value = '  +0.42%  '

# Row 15
$ws.Range("D15").Value = '3.718.42'
$ws.Range("E15").Value = '  +0.57%  '

# Row 16
$ws.Range("D16").Value = '65.897.28'
$ws.Range("E16").Value = '  +1.05%  '

# Row 17
$ws.Range("D17").Value = '7.32'
$ws.Range("E17").Value = '  +1.49%  '

# Row 18
$ws.Range("D18").Value = '3.189.15'
$ws.Range("E18").Value = '  +0.49%  '

# Row 19
$ws.Range("E19").Value = '  +0.26%  '

# Row 20
$ws.Range("D20").Value = '511.57'
$ws.Range("E20").Value = '  -0.29%  '

# Row 21
$ws.Range("D21").Value = '15.77'
$ws.Range("E21").Value = '  +5.31%  '

# Row 22
$ws.Range("D22").Value = '0.737'
$ws.Range("E22").Value = '  -0.15%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '7.96'
$ws.Range("E23").Value = '  +1.43%  '

# Row 24
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '15.11'
$ws.Range("E24").Value = '  -3.27%  '

# Row 25
$ws.Range("D25").Value = '85.34'
$ws.Range("E25").Value = '  +0.31%  '

# Row 26
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.20%  '

# Row 27
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '3.01'
$ws.Range("E27").Value = '  +3.04%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '9.21'
$ws.Range("E28").Value = '  +0.19%  '

# Row 29
$ws.Range("E29").Value = '  +1.21%  '

# Row 30
$ws.Range("D30").Value = '2.88'
$ws.Range("E30").Value = '  +3.02%  '

# Row 31
$ws.Range("D31").Value = '28.05'
$ws.Range("E31").Value = '  -0.48%  '

# Row 32
$ws.Range("D32").Value = '6.67'
$ws.Range("E32").Value = '  +5.05%  '

# Row 33
$ws.Range("E33").Value = '  +0.24%  '

# Row 34
$ws.Range("E34").Value = '  +0.13%  '

# Row 35
$ws.Range("E35").Value = '  -1.43%  '

# Row 36
$ws.Range("D36").Value = '55.57'
$ws.Range("E36").Value = '  -0.31%  '

# Row 37
$ws.Range("D37").Value = '0.0924'
$ws.Range("E37").Value = '  +3.36%  '

# Row 38
$ws.Range("D38").Value = '484.98'
$ws.Range("E38").Value = '  +0.79%  '

# Row 39
$ws.Range("D39").Value = '0.0423'
$ws.Range("E39").Value = '  +0.24%  '

# Row 40
$ws.Range("E40").Value = '  -3.37%  '

# Row 41
$ws.Range("D41").Value = '8.88'
$ws.Range("E41").Value = '  +2.65%  '

# Row 42
$ws.Range("D42").Value = '3.018.56'
$ws.Range("E42").Value = '  -3.98%  '

# Row 43
$ws.Range("D43").Value = '0.119'
$ws.Range("E43").Value = '  -1.14%  '

# Row 44
$ws.Range("D44").Value = '0.290'
$ws.Range("E44").Value = '  +0.36%  '

# Row 45
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0649'
$ws.Range("E45").Value = '  +8.13%  '

# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").Value = '  -2.58%  '

# Row 47
$ws.Range("D47").Value = '29.00'
$ws.Range("E47").Value = '  -1.24%  '

# Row 48
$ws.Range("E48").Value = '  +0.03%  '

# Row 49
$ws.Range("E49").Value = '  +0.15%  '

# Row 50
$ws.Range("D50").Value = '2.32'
$ws.Range("E50").Value = '  +0.87%  '

# Row 51
$ws.Range("D51").Value = '120.20'
$ws.Range("E51").Value = '  -1.96%  '
